$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.267.67'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.267.22'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.62'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.30'
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('E7').Value = '  -0.77%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.93'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('D14').Value = '2.619.98'
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.56'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('D16').Value = '2.269.53'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.786'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '42.173.34'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.26'
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.95'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.64'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.32'
$ws.Range('E23').Value = '  -2.51%  '
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.54'
$ws.Range('E27').Value = '  -1.57%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.02'
$ws.Range('E28').Value = '  -2.81%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.33'
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.23'
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.09'
$ws.Range('E34').Value = '  -1.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.62'
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0733'
$ws.Range('E36').Value = '  -2.12%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  -2.79%  '
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.81'
$ws.Range('E40').Value = '  -2.44%  '
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  -4.89%  '
$ws.Range('D43').Value = '1.947.91'
$ws.Range('E43').Value = '  -2.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0282'
$ws.Range('E44').Value = '  -1.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.83'
$ws.Range('E45').Value = '  -1.67%  '
$ws.Range('E46').Value = '  -2.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.76'
$ws.Range('E47').Value = '  -3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.29'
$ws.Range('E48').Value = '  +1.43%  '
$ws.Range('D49').Value = '2.492.16'
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '91.86'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.46'
$ws.Range('E51').Value = '  -2.44%  '
